# Commit: "Show credible interval in MCMCglmm analyses"
#
# 1) Rename the "Species.level.data" sheet to "Species.level.moderators"
#    (this also updates the sheet's defined-name references automatically).
# 2) Make the renamed sheet the active / selected tab (it was previously
#    the "Description" sheet that was active).
# 3) Give the "Description" sheet an explicit page setup (A4 paper,
#    portrait orientation) - it previously had no pageSetup element.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Species.level.data")
$wsData.Name = "Species.level.moderators"
$wsData.Activate()

$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.PageSetup.PaperSize = 9
$wsDesc.PageSetup.Orientation = 1
